# feat: neural networks with optimization of parameters
#
# Adds a new worksheet "neural_networks" (placed after the existing
# "Sheet1") containing an accuracy/F1 comparison table for various
# sklearn MLPClassifier solver/activation combinations, and updates the
# selection state on both sheets to match the authored workbook.

$wb = $excel.ActiveWorkbook

# --- locate the existing sheet and insert the new one right after it ---
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "neural_networks"

# --- header row (A1/B1 first; C1 "Accuracy" is written later, see below) ---
$ws2.Range("A1").Value = "metoda"
$ws2.Range("B1").Value = "F1"

# --- data rows: method name, F1, Accuracy ---
# Values in columns B/C are long decimal strings in the source workbook
# (stored as text, not numbers) -- use a leading apostrophe to force
# text entry, then strip the resulting "quote prefix" style so the cell
# keeps the workbook's default (unstyled) formatting.
$byRow = @{
    2  = @("adam-identity",   "0.64797507788161979", "0.54800000000000004")
    3  = @("adam-logistic",   "0.77551020408163263", "0.78000000000000003")
    4  = @("lbfgs-logistic",  "0.69795918367346943", "0.70399999999999996")
    5  = @("lbfgs-tanh",      "0.72332015810276684", "0.71999999999999997")
    6  = @("adam-relu",       "0.71017274472168901", "0.69799999999999995")
    7  = @("lbfgs-relu",      "0.63752276867030966", "0.60199999999999998")
    8  = @("lbfgs-identity",  "0.61056105610561062", "0.52800000000000002")
    9  = @("adam-tanh",       "0.79918032786885251", "0.80400000000000005")
    10 = @("sgd-identity",    "0.0",                 "0.496")
    11 = @("sgd-tanh",        "0.6851485148514852",  "0.68200000000000005")
    12 = @("sgd-relu",        "0.67021276595744683", "0.504")
    13 = @("sgd-logistic",    "0.0",                 "0.496")
}

# Column A (method names) was authored top-to-bottom except for row 9
# ("adam-tanh"), which was entered right after row 4 and the rows were
# later reordered -- replicate that exact entry order so the resulting
# shared-string table matches the original workbook byte-for-byte.
$colAOrder = @(2, 3, 4, 9, 5, 6, 7, 8, 10, 11, 12, 13)
foreach ($r in $colAOrder) {
    $ws2.Range("A$r").Value = $byRow[$r][0]
}

# Row 2's F1/Accuracy were entered together, then the rest of column B
# (F1) top-to-bottom, then the rest of column C (Accuracy) top-to-bottom.
$ws2.Range("B2").Value = "'" + $byRow[2][1]
$ws2.Range("C2").Value = "'" + $byRow[2][2]
for ($r = 3; $r -le 12; $r++) {
    $ws2.Range("B$r").Value = "'" + $byRow[$r][1]
}
for ($r = 3; $r -le 12; $r++) {
    $ws2.Range("C$r").Value = "'" + $byRow[$r][2]
}
$ws2.Range("B13").Value = "'" + $byRow[13][1]
$ws2.Range("C13").Value = "'" + $byRow[13][2]

# Header C1 ("Accuracy") written last.
$ws2.Range("C1").Value = "Accuracy"

# Strip the auto-applied "quote prefix" text style from the numeric-look
# text cells so they keep plain/default cell formatting.
$ws2.Range("B2:C13").ClearFormats()

# --- column widths ---
$ws2.Columns.Item(1).ColumnWidth = 28.666666666666668
$ws2.Columns.Item(2).ColumnWidth = 18.998697916666668
$ws2.Columns.Item(3).ColumnWidth = 18.998697916666668

# --- selections / active sheet ---
$null = $ws1.Range("A18").Select()
$null = $ws2.Range("A3").Select()
